$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure every written cell keeps plain-text semantics (matches original inlineStr cells)
function Set-TextCell($addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
}

Set-TextCell "D2" "30.330.44"
Set-TextCell "E2" "  -1.02%  "
Set-TextCell "D3" "1.862.06"
Set-TextCell "E3" "  -1.06%  "
Set-TextCell "D4" "0.9997"
Set-TextCell "E4" "  +0.09%  "
Set-TextCell "D5" "234.13"
Set-TextCell "E5" "  -2.21%  "
Set-TextCell "D6" "0.9998"
Set-TextCell "E6" "  +0.07%  "
Set-TextCell "D7" "0.4762"
Set-TextCell "E7" "  -0.92%  "
Set-TextCell "D8" "0.2752"
Set-TextCell "E8" "  -3.02%  "
Set-TextCell "D9" "0.06447"
Set-TextCell "E9" "  -1.56%  "
Set-TextCell "D10" "1.852.51"
Set-TextCell "E10" "  -11.74%  "
Set-TextCell "D11" "0.07430"
Set-TextCell "E11" "  -0.87%  "
Set-TextCell "D12" "16.14"
Set-TextCell "E12" "  -3.45%  "
Set-TextCell "D13" "4.998"
Set-TextCell "E13" "  -2.08%  "
Set-TextCell "D14" "86.08"
Set-TextCell "E14" "  -2.91%  "
Set-TextCell "D15" "0.6340"
Set-TextCell "E15" "  -4.90%  "
Set-TextCell "D16" "30.311.30"
Set-TextCell "E16" "  -0.93%  "
Set-TextCell "D17" "0.9997"
Set-TextCell "E17" "  +0.00%  "
Set-TextCell "D18" "232.52"
Set-TextCell "E18" "  +3.56%  "
Set-TextCell "D19" "12.83"
Set-TextCell "E19" "  -4.15%  "
Set-TextCell "D20" "0.000007398"
Set-TextCell "E20" "  -2.98%  "
Set-TextCell "B21" "WrappedliquidstakedEther2.0"
Set-TextCell "C21" "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextCell "D21" "2.096.07"
Set-TextCell "E21" "  -2.68%  "
Set-TextCell "B22" "BinanceUSD"
Set-TextCell "C22" "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextCell "D22" "0.9998"
Set-TextCell "E22" "  +0.09%  "
Set-TextCell "B23" "Uniswap"
Set-TextCell "C23" "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextCell "D23" "5.112"
Set-TextCell "E23" "  -4.36%  "
Set-TextCell "B24" "Chainlink"
Set-TextCell "C24" "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextCell "D24" "6.022"
Set-TextCell "E24" "  -3.54%  "
Set-TextCell "B25" "Cosmos"
Set-TextCell "C25" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextCell "D25" "9.300"
Set-TextCell "E25" "  -0.57%  "
Set-TextCell "B26" "Monero"
Set-TextCell "C26" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextCell "D26" "167.52"
Set-TextCell "E26" "  +0.45%  "
Set-TextCell "B27" "EthereumClassic"
Set-TextCell "C27" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextCell "D27" "17.92"
Set-TextCell "E27" "  -3.80%  "
Set-TextCell "B28" "LidoDAOToken"
Set-TextCell "C28" "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextCell "D28" "1.862"
Set-TextCell "E28" "  -5.33%  "
Set-TextCell "D29" "1.384"
Set-TextCell "E29" "  -4.97%  "
Set-TextCell "B30" "Stellar"
Set-TextCell "C30" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextCell "D30" "0.1007"
Set-TextCell "E30" "  +5.79%  "
Set-TextCell "B31" "InternetComputer(DFINITY)"
Set-TextCell "C31" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextCell "D31" "4.238"
Set-TextCell "E31" "  -2.46%  "
Set-TextCell "B32" "Filecoin"
Set-TextCell "C32" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextCell "D32" "3.921"
Set-TextCell "E32" "  -3.03%  "
Set-TextCell "B33" "Hedera"
Set-TextCell "C33" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextCell "D33" "0.04904"
Set-TextCell "E33" "  -2.77%  "
Set-TextCell "B34" "ARBITRUM"
Set-TextCell "C34" "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextCell "D34" "1.151"
Set-TextCell "E34" "  -4.91%  "
Set-TextCell "B35" "ImmutableX"
Set-TextCell "C35" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextCell "D35" "0.7273"
Set-TextCell "E35" "  -3.51%  "
Set-TextCell "B36" "Frax"
Set-TextCell "C36" "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextCell "D36" "0.9988"
Set-TextCell "E36" "  -1.47%  "
Set-TextCell "B37" "HuobiToken"
Set-TextCell "C37" "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextCell "D37" "2.693"
Set-TextCell "E37" "  -0.51%  "
Set-TextCell "B38" "VeChain"
Set-TextCell "C38" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextCell "D38" "0.01947"
Set-TextCell "E38" "  +5.92%  "
Set-TextCell "B39" "MXToken"
Set-TextCell "C39" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextCell "D39" "2.632"
Set-TextCell "E39" "  +0.20%  "
Set-TextCell "B40" "TrustWalletToken"
Set-TextCell "C40" "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextCell "D40" "0.9099"
Set-TextCell "E40" "  -0.18%  "
Set-TextCell "B41" "RenderToken"
Set-TextCell "C41" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell "D41" "1.994"
Set-TextCell "E41" "  -4.58%  "
Set-TextCell "B42" "Quant"
Set-TextCell "C42" "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextCell "D42" "105.66"
Set-TextCell "E42" "  -0.24%  "
Set-TextCell "B43" "PaxDollar"
Set-TextCell "C43" "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextCell "D43" "1.000"
Set-TextCell "E43" "  -0.57%  "
Set-TextCell "B44" "TheSandbox"
Set-TextCell "C44" "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextCell "D44" "0.4123"
Set-TextCell "E44" "  -4.10%  "
Set-TextCell "B45" "FraxShare"
Set-TextCell "C45" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextCell "D45" "5.554"
Set-TextCell "E45" "  -5.34%  "
Set-TextCell "B46" "Aptos"
Set-TextCell "C46" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextCell "D46" "7.079"
Set-TextCell "E46" "  -5.63%  "
Set-TextCell "B47" "Aave"
Set-TextCell "C47" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextCell "D47" "61.50"
Set-TextCell "E47" "  -5.06%  "
Set-TextCell "B48" "Algorand"
Set-TextCell "C48" "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextCell "D48" "0.1209"
Set-TextCell "E48" "  -6.09%  "
Set-TextCell "B49" "EnergySwap"
Set-TextCell "C49" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextCell "D49" "8.759"
Set-TextCell "E49" "  -2.63%  "
Set-TextCell "B50" "NEARProtocol"
Set-TextCell "C50" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextCell "D50" "1.404"
Set-TextCell "E50" "  -4.91%  "
Set-TextCell "D51" "33.10"
Set-TextCell "E51" "  -2.35%  "
